# Update the "Generate Report for Handback" timestamps recorded in the
# handback-status workbook.
#
# Overview!G2                -> Latest HO Xliff Generate Date (de1badf3... file)
# zh-cn!H2                   -> Correspond Handoff Datetime   (de1badf3... zh-cn xlf)
# zh-cn!K2                   -> Correspond Handback DateTime  (de1badf3... zh-cn xlf)
# de-de!H2                   -> Correspond Handoff Datetime   (de1badf3... de-de xlf)
# de-de!K2                   -> Correspond Handback DateTime  (de1badf3... de-de xlf)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-09-05 21:16:21"

$zhcn.Range("H2").Value = "2016-09-05 21:16:15"
$zhcn.Range("K2").Value = "2016-09-05 21:16:33"

$dede.Range("H2").Value = "2016-09-05 21:16:21"
$dede.Range("K2").Value = "2016-09-05 21:16:41"
